$wb = $excel.ActiveWorkbook
$wsInputs = $wb.Worksheets.Item("Inputs")
$wsDemand = $wb.Worksheets.Item("Demand_calc")

# --- Inputs sheet: insert a new row 6 for the "Week" resolution ---
# (this pushes old rows 6,7,8 -> 7,8,9, and Excel auto-updates any
# formulas elsewhere that reference Inputs!G6/G7/G8)
$wsInputs.Rows.Item(6).Insert()

$wsInputs.Range("G6").Formula = "=365/7"

# --- Demand_calc sheet: insert a new column H for the "Week" resolution ---
# (this pushes old columns H,I -> I,J)
$wsDemand.Columns.Item(8).Insert()

# Header label for the new column (added first so it lands before
# "weeks/a" in the shared-string table, matching authoring order)
$wsDemand.Range("H1").Value = "Week"

# Unit label for the new Inputs row
$wsInputs.Range("H6").Value = "weeks/a"

# New Demand_calc!H2 value (mirrors the committed formula exactly)
$wsDemand.Range("H2").Formula = "=A3/Inputs!G8"

# --- selection / active-cell bookkeeping to mirror the saved view state ---
$wsInputs.Activate()
$wsInputs.Range("G6").Select()
$wsDemand.Activate()
$wsDemand.Range("J7").Select()
